# "show & edit student grade"
# The "Point" column header is being replaced with "Fullname", and the
# active selection is moved up one row (from B3 to B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 held the shared string "Point" -> rename the column header to "Fullname"
$ws.Range("B1").Value = "Fullname"

# Move/save the active cell selection to B2 (was B3)
$ws.Range("B2").Select()
